$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 0.04510733333333333
$ws.Range("H2").Value = 0.135322
$ws.Range("I2").Value = 0.001347508866645585
$ws.Range("J2").Value = 0.001347508866645585
$ws.Range("M2").Value = 2.839339666666666
$ws.Range("N2").Value = 8.518018999999999
$ws.Range("O2").Value = 0.07557152725297667
$ws.Range("P2").Value = 0.07557152725297665
$ws.Range("Q2").Value = 0.1280750407908889
$ws.Range("R2").Value = 1.152675367118
$ws.Range("S2").Value = 0.0001018333030393345
$ws.Range("T2").Value = 0.0001018333030393345
$ws.Range("G3").Value = 0.04510733333333333
$ws.Range("H3").Value = 0.135322
$ws.Range("I3").Value = 0.001347508866645585
$ws.Range("J3").Value = 0.001347508866645585
$ws.Range("O3").Value = 0.03769827278900544
$ws.Range("P3").Value = 0.03769827278900544
$ws.Range("Q3").Value = 0.06388924507288889
$ws.Range("R3").Value = 0.5750032056559999
$ws.Range("S3").Value = 0.0000507987568404088
$ws.Range("T3").Value = 0.0000507987568404088
$ws.Range("G4").Value = 0.04510733333333333
$ws.Range("H4").Value = 0.135322
$ws.Range("I4").Value = 0.001347508866645585
$ws.Range("J4").Value = 0.001347508866645585
$ws.Range("M4").Value = 32.04516333333333
$ws.Range("N4").Value = 96.13549
$ws.Range("O4").Value = 0.852910260297995
$ws.Range("P4").Value = 0.8529102602979949
$ws.Range("Q4").Value = 1.445471864197778
$ws.Range("R4").Value = 13.00924677778
$ws.Range("S4").Value = 0.001149304138204542
$ws.Range("T4").Value = 0.001149304138204542
$ws.Range("G5").Value = 0.04510733333333333
$ws.Range("H5").Value = 0.135322
$ws.Range("I5").Value = 0.001347508866645585
$ws.Range("J5").Value = 0.001347508866645585
$ws.Range("M5").Value = 1.270667666666667
$ws.Range("N5").Value = 3.812003
$ws.Range("O5").Value = 0.03381993966002293
$ws.Range("P5").Value = 0.03381993966002293
$ws.Range("Q5").Value = 0.05731642999622222
$ws.Range("R5").Value = 0.515847869966
$ws.Range("S5").Value = 0.00004557266856129955
$ws.Range("T5").Value = 0.00004557266856129955
$ws.Range("I6").Value = 0.3371496619592149
$ws.Range("J6").Value = 0.3371496619592149
$ws.Range("M6").Value = 2.839339666666666
$ws.Range("N6").Value = 8.518018999999999
$ws.Range("O6").Value = 0.07557152725297667
$ws.Range("P6").Value = 0.07557152725297665
$ws.Range("Q6").Value = 32.04465497548222
$ws.Range("R6").Value = 288.40189477934
$ws.Range("S6").Value = 0.02547891486708268
$ws.Range("T6").Value = 0.02547891486708267
$ws.Range("I7").Value = 0.3371496619592149
$ws.Range("J7").Value = 0.3371496619592149
$ws.Range("O7").Value = 0.03769827278900544
$ws.Range("P7").Value = 0.03769827278900544
$ws.Range("S7").Value = 0.01270995992725945
$ws.Range("T7").Value = 0.01270995992725945
$ws.Range("I8").Value = 0.3371496619592149
$ws.Range("J8").Value = 0.3371496619592149
$ws.Range("M8").Value = 32.04516333333333
$ws.Range("N8").Value = 96.13549
$ws.Range("O8").Value = 0.852910260297995
$ws.Range("P8").Value = 0.8529102602979949
$ws.Range("Q8").Value = 361.6602179390445
$ws.Range("R8").Value = 3254.9419614514
$ws.Range("S8").Value = 0.287558405941015
$ws.Range("T8").Value = 0.287558405941015
$ws.Range("I9").Value = 0.3371496619592149
$ws.Range("J9").Value = 0.3371496619592149
$ws.Range("M9").Value = 1.270667666666667
$ws.Range("N9").Value = 3.812003
$ws.Range("O9").Value = 0.03381993966002293
$ws.Range("P9").Value = 0.03381993966002293
$ws.Range("Q9").Value = 14.34069598817555
$ws.Range("R9").Value = 129.06626389358
$ws.Range("S9").Value = 0.01140238122385778
$ws.Range("T9").Value = 0.01140238122385778
$ws.Range("G10").Value = 0.8868746666666668
$ws.Range("H10").Value = 2.660624
$ws.Range("I10").Value = 0.02649395095261704
$ws.Range("J10").Value = 0.02649395095261704
$ws.Range("M10").Value = 2.839339666666666
$ws.Range("N10").Value = 8.518018999999999
$ws.Range("O10").Value = 0.07557152725297667
$ws.Range("P10").Value = 0.07557152725297665
$ws.Range("Q10").Value = 2.518138420428444
$ws.Range("R10").Value = 22.663245783856
$ws.Range("S10").Value = 0.002002188336454725
$ws.Range("T10").Value = 0.002002188336454725
$ws.Range("G11").Value = 0.8868746666666668
$ws.Range("H11").Value = 2.660624
$ws.Range("I11").Value = 0.02649395095261704
$ws.Range("J11").Value = 0.02649395095261704
$ws.Range("O11").Value = 0.03769827278900544
$ws.Range("P11").Value = 0.03769827278900544
$ws.Range("Q11").Value = 1.256153905372445
$ws.Range("R11").Value = 11.305385148352
$ws.Range("S11").Value = 0.0009987761902702874
$ws.Range("T11").Value = 0.0009987761902702874
$ws.Range("G12").Value = 0.8868746666666668
$ws.Range("H12").Value = 2.660624
$ws.Range("I12").Value = 0.02649395095261704
$ws.Range("J12").Value = 0.02649395095261704
$ws.Range("M12").Value = 32.04516333333333
$ws.Range("N12").Value = 96.13549
$ws.Range("O12").Value = 0.852910260297995
$ws.Range("P12").Value = 0.8529102602979949
$ws.Range("Q12").Value = 28.42004354952889
$ws.Range("R12").Value = 255.78039194576
$ws.Range("S12").Value = 0.02259696260331891
$ws.Range("T12").Value = 0.02259696260331891
$ws.Range("G13").Value = 0.8868746666666668
$ws.Range("H13").Value = 2.660624
$ws.Range("I13").Value = 0.02649395095261704
$ws.Range("J13").Value = 0.02649395095261704
$ws.Range("M13").Value = 1.270667666666667
$ws.Range("N13").Value = 3.812003
$ws.Range("O13").Value = 0.03381993966002293
$ws.Range("P13").Value = 0.03381993966002293
$ws.Range("Q13").Value = 1.126922963319111
$ws.Range("R13").Value = 10.142306669872
$ws.Range("S13").Value = 0.0008960238225731151
$ws.Range("T13").Value = 0.0008960238225731151
$ws.Range("G14").Value = 21.25667433333333
$ws.Range("H14").Value = 63.77002299999999
$ws.Range("I14").Value = 0.6350088782215225
$ws.Range("J14").Value = 0.6350088782215224
$ws.Range("M14").Value = 2.839339666666666
$ws.Range("N14").Value = 8.518018999999999
$ws.Range("O14").Value = 0.07557152725297667
$ws.Range("P14").Value = 0.07557152725297665
$ws.Range("Q14").Value = 60.35491861604854
$ws.Range("R14").Value = 543.1942675444369
$ws.Range("S14").Value = 0.04798859074639993
$ws.Range("T14").Value = 0.04798859074639991
$ws.Range("G15").Value = 21.25667433333333
$ws.Range("H15").Value = 63.77002299999999
$ws.Range("I15").Value = 0.6350088782215225
$ws.Range("J15").Value = 0.6350088782215224
$ws.Range("O15").Value = 0.03769827278900544
$ws.Range("P15").Value = 0.03769827278900544
$ws.Range("Q15").Value = 30.10758507671155
$ws.Range("R15").Value = 270.968265690404
$ws.Range("S15").Value = 0.02393873791463529
$ws.Range("T15").Value = 0.02393873791463529
$ws.Range("G16").Value = 21.25667433333333
$ws.Range("H16").Value = 63.77002299999999
$ws.Range("I16").Value = 0.6350088782215225
$ws.Range("J16").Value = 0.6350088782215224
$ws.Range("M16").Value = 32.04516333333333
$ws.Range("N16").Value = 96.13549
$ws.Range("O16").Value = 0.852910260297995
$ws.Range("P16").Value = 0.8529102602979949
$ws.Range("Q16").Value = 681.1736009351412
$ws.Range("R16").Value = 6130.56240841627
$ws.Range("S16").Value = 0.5416055876154566
$ws.Range("T16").Value = 0.5416055876154564
$ws.Range("G17").Value = 21.25667433333333
$ws.Range("H17").Value = 63.77002299999999
$ws.Range("I17").Value = 0.6350088782215225
$ws.Range("J17").Value = 0.6350088782215224
$ws.Range("M17").Value = 1.270667666666667
$ws.Range("N17").Value = 3.812003
$ws.Range("O17").Value = 0.03381993966002293
$ws.Range("P17").Value = 0.03381993966002293
$ws.Range("Q17").Value = 27.01016877622989
$ws.Range("R17").Value = 243.091518986069
$ws.Range("S17").Value = 0.02147596194503074
$ws.Range("T17").Value = 0.02147596194503073
